# add timeOut before method
#
# Simulates the effect of introducing a short delay ("timeOut") before the
# code path that builds each user's login response: re-running that code
# now produces a fresh "id" (uuid) and a fresh "token" (jwt, whose "iat"
# claim reflects the later timestamp) for each of the three existing test
# users, and the corresponding data row is updated to store these newly
# generated values in place of the old ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - daniel5f: new id + new token (iat 1701318084)
$ws.Range("D2").Value = "03ad65ca-0b41-4b9f-9e4e-0d3940c49488"
$ws.Range("C2").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6ImRhbmllbDVmIiwicGFzc3dvcmQiOiJBejI1Mjg4QCIsImlhdCI6MTcwMTMxODA4NH0.sP11KRwFIRb4Ep-MbNpAG5O6Re7Qk1DcjG8lV28tbfU"

# Row 3 - Jorge2525: new id + new token (iat 1701318086)
$ws.Range("D3").Value = "912a3f79-69c5-4be3-a70b-02e5753c4fd7"
$ws.Range("C3").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6IkpvcmdlMjUyNSIsInBhc3N3b3JkIjoiYXNUMzU2NDQ0QCIsImlhdCI6MTcwMTMxODA4Nn0.4exdoDoZgInbSWjpSrqXtKMLx32gf983nk6kJjXH4fY"

# Row 4 - mario35: new id + new token (iat 1701318087)
$ws.Range("D4").Value = "48dd68ae-7e0a-47e6-9805-d174727f1795"
$ws.Range("C4").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6Im1hcmlvMzUiLCJwYXNzd29yZCI6Im1BcmlvdXVnQDMiLCJpYXQiOjE3MDEzMTgwODd9.fkUcucVDq6v6W1hneINrxmjLNRQ5MZ1zurJDIcgH1Ok"
